$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert 3 missing-quarter blocks (Q2-2016, Q4-2015, Q1-2016), each a
# pair of rows (F&O, CDS), at the positions identified from the diff.
# Inserted from bottom to top so earlier row numbers stay valid.
# ------------------------------------------------------------------

# Insert Q1-2016 block just above old row 16 (Q2-2018 F&O) -> becomes rows 20-21
$ws.Rows("16:17").Insert()

# Insert Q4-2015 block just above old row 14 (Q2-2017 F&O) -> becomes rows 16-17
$ws.Rows("14:15").Insert()

# Insert Q2-2016 block just above old row 8 (Q4-2017 F&O) -> becomes rows 8-9
$ws.Rows("8:9").Insert()

# Newly-inserted rows do not inherit the bordered/bold style used by column A
# in the data rows (style of A2); re-apply it explicitly.
$ws.Range("A8").Style = $ws.Range("A2").Style
$ws.Range("A9").Style = $ws.Range("A2").Style
$ws.Range("A16").Style = $ws.Range("A2").Style
$ws.Range("A17").Style = $ws.Range("A2").Style
$ws.Range("A20").Style = $ws.Range("A2").Style
$ws.Range("A21").Style = $ws.Range("A2").Style

# --- Populate data for newly inserted rows ---
# Row 8: Q2-2016 F&O
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "ICE CLEAR EUROPE"
$ws.Range("C8").Value = "Q2-2016"
$ws.Range("D8").Value = "F&O"
$ws.Range("E8").Value = 90073860
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 1666366420
$ws.Range("I8").Value = 1728696812
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("L8").Value = 3332732840
$ws.Range("N8").Value = 0
$ws.Range("Q8").Value = 0
$ws.Range("R8").Value = 34405787647
$ws.Range("U8").Value = 72
$ws.Range("V8").Value = 0
$ws.Range("W8").Value = 0
$ws.Range("X8").Value = "Not available"
$ws.Range("Y8").Value = 64
$ws.Range("Z8").Value = 0
$ws.Range("AA8").Value = 0
$ws.Range("AB8").Value = 0
$ws.Range("AC8").Value = 0

# Row 9: Q2-2016 CDS
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "ICE CLEAR EUROPE"
$ws.Range("C9").Value = "Q2-2016"
$ws.Range("D9").Value = "CDS"
$ws.Range("E9").Value = 22508328
$ws.Range("F9").Value = 22508328
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 973277475
$ws.Range("I9").Value = 993159372
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 973277475
$ws.Range("N9").Value = 0
$ws.Range("O9").Value = 160782997
$ws.Range("P9").Value = 62273368
$ws.Range("Q9").Value = 284724655
$ws.Range("R9").Value = 7016438114
$ws.Range("U9").Value = 22
$ws.Range("V9").Value = 0
$ws.Range("W9").Value = 0
$ws.Range("X9").Value = "Not available"
$ws.Range("Y9").Value = 8
$ws.Range("Z9").Value = 10657801843
$ws.Range("AA9").Value = 0
$ws.Range("AC9").Value = 0

# Row 16: Q4-2015 F&O
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "ICE CLEAR EUROPE"
$ws.Range("C16").Value = "Q4-2015"
$ws.Range("D16").Value = "F&O"
$ws.Range("E16").Value = 91852668
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 1707165908
$ws.Range("I16").Value = 1755379393
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 3414331817
$ws.Range("N16").Value = 0
$ws.Range("Q16").Value = 0
$ws.Range("R16").Value = 37121934837
$ws.Range("U16").Value = 73
$ws.Range("V16").Value = 0
$ws.Range("W16").Value = 0
$ws.Range("X16").Value = "Not available"
$ws.Range("Y16").Value = 62
$ws.Range("Z16").Value = 0
$ws.Range("AA16").Value = 0
$ws.Range("AB16").Value = 0
$ws.Range("AC16").Value = 0

# Row 17: Q4-2015 CDS
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "ICE CLEAR EUROPE"
$ws.Range("C17").Value = "Q4-2015"
$ws.Range("D17").Value = "CDS"
$ws.Range("E17").Value = 22956841
$ws.Range("F17").Value = 22956841
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 843166303
$ws.Range("I17").Value = 855078752
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 843166303
$ws.Range("N17").Value = 0
$ws.Range("Q17").Value = 270072534
$ws.Range("R17").Value = 6053614803
$ws.Range("U17").Value = 22
$ws.Range("V17").Value = 0
$ws.Range("W17").Value = 0
$ws.Range("X17").Value = "Not available"
$ws.Range("Y17").Value = 5
$ws.Range("Z17").Value = 8669177383
$ws.Range("AA17").Value = 0
$ws.Range("AC17").Value = 0

# Row 20: Q1-2016 F&O
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "ICE CLEAR EUROPE"
$ws.Range("C20").Value = "Q1-2016"
$ws.Range("D20").Value = "F&O"
$ws.Range("E20").Value = 87834870
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 1624945103
$ws.Range("I20").Value = 1675297057
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 3249890206
$ws.Range("N20").Value = 0
$ws.Range("Q20").Value = 0
$ws.Range("R20").Value = 35264934167
$ws.Range("U20").Value = 73
$ws.Range("V20").Value = 0
$ws.Range("W20").Value = 0
$ws.Range("X20").Value = "Not available"
$ws.Range("Y20").Value = 62
$ws.Range("Z20").Value = 0
$ws.Range("AA20").Value = 0
$ws.Range("AB20").Value = 0
$ws.Range("AC20").Value = 0

# Row 21: Q1-2016 CDS
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "ICE CLEAR EUROPE"
$ws.Range("C21").Value = "Q1-2016"
$ws.Range("D21").Value = "CDS"
$ws.Range("E21").Value = 21968365
$ws.Range("F21").Value = 21968365
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 1125918062
$ws.Range("I21").Value = 1139429731
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 1125918062
$ws.Range("N21").Value = 0
$ws.Range("Q21").Value = 262851164
$ws.Range("R21").Value = 6382175808
$ws.Range("U21").Value = 22
$ws.Range("V21").Value = 0
$ws.Range("W21").Value = 0
$ws.Range("X21").Value = "Not available"
$ws.Range("Y21").Value = 7
$ws.Range("Z21").Value = 11076652084
$ws.Range("AA21").Value = 0
$ws.Range("AC21").Value = 0

# ------------------------------------------------------------------
# Renumber column A (row index) sequentially 0..21 for all data rows,
# reflecting the new row order after insertion.
# ------------------------------------------------------------------
for ($r = 2; $r -le 23; $r++) {
    $ws.Range("A$r").Value = $r - 2
}

